$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Developer name (student name) in the header
$ws.Range("C3").Value = 'Sukhtab Singh Warya'

# Test plan rows for the Client class __init__ / accessor / __str__ tests
$ws.Range("E7").Value = 'None'
$ws.Range("F7").Value = 'client_number=1010, first_name="Susan", last_name="Clark", email_address="susanclark@pixell.com"'
$ws.Range("G7").Value = 'client_number should be 1010, first_name should be "Susan", last_name should be "Clark", email_address should be "susanclark@pixell.com"'
$ws.Range("E8").Value = 'None'
$ws.Range("F8").Value = 'client_number="not_integer", first_name="Susan", last_name="Clark", email_address="susanclark@pixell.com"'
$ws.Range("G8").Value = 'ValueError with message "Client number must be an integer."'
$ws.Range("E9").Value = 'None'
$ws.Range("F9").Value = 'client_number=1010, first_name=" ", last_name="Clark", email_address="susanclark@pixell.com"'
$ws.Range("G9").Value = 'ValueError with message "First name cannot be blank."'
$ws.Range("E10").Value = 'None'
$ws.Range("F10").Value = 'client_number=1010, first_name="Susan", last_name=" ", email_address="susanclark@pixell.com"'
$ws.Range("G10").Value = 'ValueError with message "Last name cannot be blank."'
$ws.Range("E11").Value = 'None'
$ws.Range("F11").Value = 'client_number=1010, first_name="Susan", last_name="Clark", email_address="invalid-email"'
$ws.Range("G11").Value = 'email_address should be "email@pixell-river.com"'
$ws.Range("E12").Value = 'Client instance is created with valid attributes.'
$ws.Range("F12").Value = 'Client instance created with client_number=1010'
$ws.Range("G12").Value = 'client_number should be 1010'
$ws.Range("E13").Value = 'Client instance is created with valid attributes.'
$ws.Range("F13").Value = 'Client instance created with first_name="Susan"'
$ws.Range("G13").Value = 'first_name should be "Susan"'
$ws.Range("E14").Value = 'Client instance is created with valid attributes.'
$ws.Range("F14").Value = 'Client instance created with last_name="Clark"'
$ws.Range("G14").Value = 'last_name should be "Clark"'
$ws.Range("E15").Value = 'Client instance is created with valid attributes.'
$ws.Range("F15").Value = 'Client instance created with email_address="susanclark@pixell.com"'
$ws.Range("G15").Value = 'email_address should be "susanclark@pixell.com"'
$ws.Range("E16").Value = 'Client instance is created with valid attributes.'
$ws.Range("F16").Value = 'Client instance created with client_number=1010, first_name="Susan", last_name="Clark", email_address="susanclark@pixell.com"'
$ws.Range("G16").Value = 'The string should be "Clark, Susan [1010] - susanclark@pixell.com\n"'

# The placeholder cells used a bold font; once real data is entered the
# template switches these to a regular (non-bold) weight.
$ws.Range("C3").Font.Bold = $false
$ws.Range("E7:G16").Font.Bold = $false

# Keep the student-name / table-intro cell selected, matching the saved view
$ws.Range("C3:D3").Select()
